# Updates stock-qty/value figures (and dependent subtotal/grand-total cells)
# in the Companywise Stock Report to match the revised inventory counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 139
$ws.Range("G6").Value = 4153.32
$ws.Range("F9").Value = 66
$ws.Range("G9").Value = 1951.62
$ws.Range("B10").Value = 35804.45
$ws.Range("F90").Value = 66
$ws.Range("G90").Value = 8906.700000000001
$ws.Range("B114").Value = 231150.97
$ws.Range("B136").Value = 48654
$ws.Range("E136").Value = 38.26
$ws.Range("F136").Value = -1
$ws.Range("G136").Value = -32.02
$ws.Range("B137").Value = 63902
$ws.Range("E137").Value = 34.04
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("B163").Value = 64329
$ws.Range("E163").Value = 128.32
$ws.Range("F163").Value = 2
$ws.Range("G163").Value = 241.38
$ws.Range("B164").Value = 57552
$ws.Range("E164").Value = 136.86
$ws.Range("F164").Value = -5
$ws.Range("G164").Value = -603.45
$ws.Range("B292").Value = 63520
$ws.Range("E292").Value = 153.4
$ws.Range("F292").Value = 69
$ws.Range("G292").Value = 9955.32
$ws.Range("B293").Value = 55373
$ws.Range("E293").Value = 163.62
$ws.Range("F293").Value = -94
$ws.Range("G293").Value = -13562.32
$ws.Range("F296").Value = 79
$ws.Range("G296").Value = 11334.92
$ws.Range("F301").Value = 52
$ws.Range("G301").Value = 2668.64
$ws.Range("F324").Value = 25
$ws.Range("G324").Value = 1478.25
$ws.Range("F328").Value = 178
$ws.Range("G328").Value = 3743.34
$ws.Range("B339").Value = 251997.98
$ws.Range("F358").Value = 3
$ws.Range("G358").Value = 963.03
$ws.Range("B361").Value = 11235.8
$ws.Range("F386").Value = 14
$ws.Range("G386").Value = 1451.66
$ws.Range("F387").Value = 121
$ws.Range("G387").Value = 2565.2
$ws.Range("B395").Value = 224527.89
$ws.Range("B420").Value = 47097
$ws.Range("D420").Value = 112.28
$ws.Range("E420").Value = 134.16
$ws.Range("F420").Value = 15
$ws.Range("G420").Value = 1684.2
$ws.Range("B421").Value = 58047
$ws.Range("D421").Value = 105.54
$ws.Range("E421").Value = 126.1
$ws.Range("F421").Value = 41
$ws.Range("G421").Value = 4327.14
$ws.Range("F426").Value = 30
$ws.Range("G426").Value = 2898
$ws.Range("F427").Value = 88
$ws.Range("G427").Value = 3274.48
$ws.Range("B430").Value = 37763.6
$ws.Range("F439").Value = 92
$ws.Range("G439").Value = 2180.4
$ws.Range("B448").Value = 35356.41
$ws.Range("F456").Value = 39
$ws.Range("G456").Value = 2675.01
$ws.Range("B460").Value = 41973.52
$ws.Range("F484").Value = 381
$ws.Range("G484").Value = 2472.69
$ws.Range("F490").Value = 175
$ws.Range("G490").Value = 2577.75
$ws.Range("B492").Value = -13964.73
$ws.Range("F545").Value = 10
$ws.Range("G545").Value = 2454.7
$ws.Range("B546").Value = 5945.08
$ws.Range("F555").Value = 65
$ws.Range("G555").Value = 2121.6
$ws.Range("B556").Value = 5254.37
$ws.Range("B568").Value = 53319
$ws.Range("E568").Value = 310.64
$ws.Range("F568").Value = -6
$ws.Range("G568").Value = -1643.52
$ws.Range("B569").Value = 64810
$ws.Range("E569").Value = 291.22
$ws.Range("F569").Value = 5
$ws.Range("G569").Value = 1369.6
$ws.Range("F582").Value = 54
$ws.Range("G582").Value = 5970.78
$ws.Range("B586").Value = 18475.18
$ws.Range("B600").Value = 64830
$ws.Range("E600").Value = 34.9
$ws.Range("F600").Value = 111
$ws.Range("G600").Value = 3644.13
$ws.Range("B601").Value = 60022
$ws.Range("E601").Value = 37.22
$ws.Range("F601").Value = -113
$ws.Range("G601").Value = -3709.79
$ws.Range("F615").Value = 40
$ws.Range("G615").Value = 7018.8
$ws.Range("B618").Value = 35843.59
$ws.Range("F703").Value = 37
$ws.Range("G703").Value = 1770.82
$ws.Range("F705").Value = 77
$ws.Range("G705").Value = 11021.01
$ws.Range("B709").Value = 63150
$ws.Range("D709").Value = 75.68000000000001
$ws.Range("E709").Value = 80.45
$ws.Range("F709").Value = 21
$ws.Range("G709").Value = 1589.28
$ws.Range("B710").Value = 61428
$ws.Range("D710").Value = 69.16
$ws.Range("E710").Value = 73.52
$ws.Range("F710").Value = 1
$ws.Range("G710").Value = 69.16
$ws.Range("F714").Value = 60
$ws.Range("G714").Value = 4173.6
$ws.Range("F719").Value = 73
$ws.Range("G719").Value = 8811.83
$ws.Range("B720").Value = 58742.08
$ws.Range("F739").Value = 46
$ws.Range("G739").Value = 1520.76
$ws.Range("F745").Value = 71
$ws.Range("G745").Value = 17167.8
$ws.Range("F746").Value = 86
$ws.Range("G746").Value = 4901.14
$ws.Range("B747").Value = 51678.57
$ws.Range("F772").Value = 2631
$ws.Range("G772").Value = 429142.41
$ws.Range("F774").Value = 539
$ws.Range("G774").Value = 152466.93
$ws.Range("F775").Value = 424
$ws.Range("G775").Value = 61331.6
$ws.Range("B779").Value = 652504.39
$ws.Range("F782").Value = 37
$ws.Range("G782").Value = 5402.37
$ws.Range("B796").Value = 62025.17
$ws.Range("B797").Value = 2438909.83
$ws.Range("B798").Value = 2438909.83
